# Update the Kruskal-Wallis group-significance results (H statistic, p-value,
# q-value) for the "skin" and "swab" comparisons with the refreshed values,
# and fix the Group 1 / Group 2 pairing labels that had drifted out of sync
# with their statistics for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"5.824242424242443"
$ws.Range("D2").Value = [double]"0.01580677460730209"
$ws.Range("E2").Value = [double]"0.1042698707275791"
$ws.Range("C3").Value = [double]"5.338954703832769"
$ws.Range("D3").Value = [double]"0.02085397414551581"
$ws.Range("E3").Value = [double]"0.1042698707275791"
$ws.Range("A4").Value = 'H (n=27)'
$ws.Range("B4").Value = 'J (n=27)'
$ws.Range("C4").Value = [double]"4.062027684249927"
$ws.Range("D4").Value = [double]"0.04385780613755715"
$ws.Range("E4").Value = [double]"0.1461926871251905"
$ws.Range("A5").Value = 'A (n=5)'
$ws.Range("B5").Value = 'H (n=27)'
$ws.Range("C5").Value = [double]"3.394612794612783"
$ws.Range("D5").Value = [double]"0.06540971903742436"
$ws.Range("A7").Value = 'H (n=27)'
$ws.Range("B7").Value = 'S (n=35)'
$ws.Range("C7").Value = [double]"1.762106324011086"
$ws.Range("D7").Value = [double]"0.1843620084122517"
$ws.Range("E7").Value = [double]"0.2633742977317882"
$ws.Range("A8").Value = 'S (n=35)'
$ws.Range("B8").Value = 'T (n=4)'
$ws.Range("C8").Value = [double]"1.928571428571445"
$ws.Range("D8").Value = [double]"0.1649148225532956"
$ws.Range("E8").Value = [double]"0.2633742977317882"
$ws.Range("C9").Value = [double]"1.125"
$ws.Range("D9").Value = [double]"0.2888443663464818"
$ws.Range("E9").Value = [double]"0.3610554579331022"
$ws.Range("A10").Value = 'J (n=27)'
$ws.Range("B10").Value = 'S (n=35)'
$ws.Range("C10").Value = [double]"0.4741244646006635"
$ws.Range("D10").Value = [double]"0.4910958161831776"
$ws.Range("E10").Value = [double]"0.5456620179813084"
$ws.Range("A11").Value = 'A (n=5)'
$ws.Range("B11").Value = 'T (n=4)'
$ws.Range("C11").Value = [double]"0.240000000000002"
$ws.Range("D11").Value = [double]"0.6242061147664044"
$ws.Range("E11").Value = [double]"0.6242061147664044"
$ws.Range("C13").Value = [double]"26.67432950191571"
$ws.Range("D13").Value = [double]"2.407995043080524e-07"
$ws.Range("E13").Value = [double]"9.540333889063698e-07"
$ws.Range("C14").Value = [double]"26.34058956916099"
$ws.Range("D14").Value = [double]"2.862100166719109e-07"
$ws.Range("E14").Value = [double]"9.540333889063698e-07"
$ws.Range("C15").Value = [double]"15.55555555555557"
$ws.Range("D15").Value = [double]"8.011587656029265e-05"
$ws.Range("E15").Value = [double]"0.0002002896914007316"
$ws.Range("C16").Value = [double]"12.87982387982388"
$ws.Range("D16").Value = [double]"0.0003321434948157462"
$ws.Range("E16").Value = [double]"0.0006642869896314924"
$ws.Range("C17").Value = [double]"10.76565656565656"
$ws.Range("D17").Value = [double]"0.001034008775208132"
$ws.Range("E17").Value = [double]"0.001723347958680221"
$ws.Range("C18").Value = [double]"10.38173018753781"
$ws.Range("D18").Value = [double]"0.001272683754682938"
$ws.Range("E18").Value = [double]"0.001818119649547055"
$ws.Range("C19").Value = [double]"7.239151193633972"
$ws.Range("D19").Value = [double]"0.007133069373914006"
$ws.Range("E19").Value = [double]"0.008916336717392507"
$ws.Range("C20").Value = [double]"2.725208349262175"
$ws.Range("D20").Value = [double]"0.09877523127984016"
$ws.Range("E20").Value = [double]"0.1097502569776002"
$ws.Range("C21").Value = [double]"0.1593257806023587"
$ws.Range("D21").Value = [double]"0.6897780140551726"
$ws.Range("E21").Value = [double]"0.6897780140551726"
